# tz: chi: fixed some issues on forms
#
# Survey sheet: remove the "Was the child's urine preserved?" (u_urine_conserve)
# question and its dependent "end note" (u_end_note), along with the
# 'relevant' conditions on the urine-result questions that referenced it.
# Also bump the repeat-group name / form version from v1.2 to v1.3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Delete the bottom row first (u_end_note, row 18) so the row 12 index
# below is still valid when we delete it next.
$ws.Rows.Item(18).EntireRow.Delete()
# Delete u_urine_conserve (row 12)
$ws.Rows.Item(12).EntireRow.Delete()

# After the two deletions, the urine-result questions (previously rows
# 13-17) are now rows 12-16; clear their 'relevant' (column H) formulas
# which referenced the now-deleted u_urine_conserve question.
$ws.Range("H12:H16").ClearContents()

# Rename the repeat group from u_1_2 to u_1_3
$ws.Range("B8").Value2 = "u_1_3"

# Update the active selection to match the post-edit cursor position
$ws.Range("A17:XFD17").Select()

# settings sheet: bump form_title / form_id from V1.2 to V1.3
$wsSettings = $wb.Worksheets.Item("settings")
$wsSettings.Range("A2").Value2 = "(2023 Nov) - 4. SCH/STH - Urine Filtration V1.3"
$wsSettings.Range("B2").Value2 = "tz_sch_sth_impact_202311_4_urine_filtration_v1_3"
